$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.875.48'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.838.85'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.19'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.68'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.06%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.105.19'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.46'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.843.18'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.851.05'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.81'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.35'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.19'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.03'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.78'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.78%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.71%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0553'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.67%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.26%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.43'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +11.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.694'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '91.07'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.341.15'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.50%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.78'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.27'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.019.75'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0665'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.35'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +17.60%  '
